$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The target cells hold text-like values (some of which, e.g. "0.17",
# "-0.01", look numeric). Temporarily force Text format so Excel stores
# them as strings instead of silently coercing them to numbers, then
# clear the formatting again afterwards so the cells end up with no
# explicit style, exactly like the rest of the original workbook.
$valueRange = $ws.Range("B2:D4")
$valueRange.NumberFormat = "@"

# Set values in column-major order (B2,B3,B4,C2,C3,C4,D2,D3,D4) to match
# the shared-string table ordering of the original workbook.

# Column B: A
$ws.Range("B2").Value = "0.17"
$ws.Range("B3").Value = "-0.01"
$ws.Range("B4").Value = "-0.09"

# Column C: FFR
$ws.Range("C2").Value = "44.29***"
$ws.Range("C3").Value = "2.21***"
$ws.Range("C4").Value = "0.98"

# Column D: LF
$ws.Range("D2").Value = "-0.89"
$ws.Range("D3").Value = "0.46***"
$ws.Range("D4").Value = "0.82*"

$valueRange.ClearFormats()
